$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 125 (shifts existing rows 125-177 down to 126-178)
$ws.Rows(125).Insert()

# Populate the newly inserted row 125 with the new record's data
$ws.Range("A125").Value = 9
$ws.Range("B125").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C125").Value = "Metropolitana"
$ws.Range("D125").Value = (Get-Date -Year 2023 -Month 8 -Day 3 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E125").Value = 13
$ws.Range("F125").Value = 100112022
$ws.Range("G125").Value = "Arveja Verde"
$ws.Range("H125").Value = "Perfection"
$ws.Range("I125").Value = "Primera"
$ws.Range("J125").Value = 52
$ws.Range("K125").Value = 24000
$ws.Range("L125").Value = 25000
$ws.Range("M125").Value = 24500
$ws.Range("N125").Value = "`$/malla 25 kilos"
$ws.Range("O125").Value = "Provincia de Huasco"
$ws.Range("P125").Value = 980
$ws.Range("Q125").Value = 25
$ws.Range("R125").Value = "Hortaliza"

# Apply the same number format as the other date cells in column D
$ws.Range("D125").NumberFormat = $ws.Range("D126").NumberFormat()
